$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 6807
$ws.Range("F3").Value = 822
$ws.Range("F5").Value = 145
$ws.Range("F7").Value = 738
$ws.Range("F8").Value = 738
$ws.Range("F12").Value = 1110
$ws.Range("F13").Value = 871
$ws.Range("F17").Value = 1352
$ws.Range("F20").Value = 540
$ws.Range("F21").Value = 6
$ws.Range("F22").Value = 573
$ws.Range("F26").Value = 1067
$ws.Range("F28").Value = 727
$ws.Range("F29").Value = 547
$ws.Range("F30").Value = 465
$ws.Range("F31").Value = 459
$ws.Range("F34").Value = 1131
$ws.Range("F35").Value = 268
$ws.Range("F38").Value = 1266
$ws.Range("F41").Value = 3878

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 750
$ws.Range("F9").Value = 162
$ws.Range("F16").Value = 338
$ws.Range("F17").Value = 4130
$ws.Range("F23").Value = 241
$ws.Range("F25").Value = 112

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 1257
$ws.Range("F5").Value = 1638

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1257
$ws.Range("F4").Value = 1638
$ws.Range("F8").Value = 6807
$ws.Range("F10").Value = 822
$ws.Range("F11").Value = 750
$ws.Range("F12").Value = 145
$ws.Range("F14").Value = 738
$ws.Range("F15").Value = 738
$ws.Range("F18").Value = 1110
$ws.Range("F19").Value = 871
$ws.Range("F21").Value = 162
$ws.Range("F22").Value = 162
$ws.Range("F25").Value = 1352
$ws.Range("F28").Value = 540
$ws.Range("F29").Value = 6
$ws.Range("F30").Value = 573
$ws.Range("F32").Value = 338
$ws.Range("F34").Value = 1067
$ws.Range("F36").Value = 727
$ws.Range("F37").Value = 547
$ws.Range("F38").Value = 465
$ws.Range("F39").Value = 459
$ws.Range("F43").Value = 1131
$ws.Range("F44").Value = 268
$ws.Range("F49").Value = 1267
$ws.Range("F51").Value = 3878
